# Update auto-increase logic to support new hire only and all eligible scenarios
# This applies the recalculated projection values for rows 2-6 (columns C:Q)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("C2").Value = 102
$ws.Range("D2").Value = 88
$ws.Range("E2").Value = 0.8627450980392157
$ws.Range("F2").Value = 0.8627450980392157
$ws.Range("G2").Value = 0.09833450573160141
$ws.Range("H2").Value = 0.08483761278804831
$ws.Range("I2").Value = 472190.6408301356
$ws.Range("J2").Value = 173481.3206160678
$ws.Range("L2").Value = 173481.3206160678
$ws.Range("M2").Value = 645671.9614462035
$ws.Range("N2").Value = 10262129.6088
$ws.Range("O2").Value = 9854388.678699998
$ws.Range("P2").Value = 0.01690500190791819
$ws.Range("Q2").Value = 0.01760447312079775

# Row 3
$ws.Range("C3").Value = 103
$ws.Range("D3").Value = 89
$ws.Range("E3").Value = 0.8640776699029126
$ws.Range("F3").Value = 0.8640776699029126
$ws.Range("G3").Value = 0.09614705683641531
$ws.Range("H3").Value = 0.08307852483923264
$ws.Range("I3").Value = 477788.110972513
$ws.Range("J3").Value = 173681.6356932865
$ws.Range("L3").Value = 173681.6356932865
$ws.Range("M3").Value = 651469.7466657992
$ws.Range("N3").Value = 10483190.267664
$ws.Range("O3").Value = 10075817.109661
$ws.Range("P3").Value = 0.01656763172838877
$ws.Range("Q3").Value = 0.01723747402349683

# Row 4
$ws.Range("C4").Value = 103
$ws.Range("D4").Value = 89
$ws.Range("E4").Value = 0.8640776699029126
$ws.Range("F4").Value = 0.8557692307692307
$ws.Range("G4").Value = 0.0970994315717198
$ws.Range("H4").Value = 0.08319085971041408
$ws.Range("I4").Value = 509184.1566766572
$ws.Range("J4").Value = 182293.2859515694
$ws.Range("L4").Value = 182293.2859515694
$ws.Range("M4").Value = 691477.4426282267
$ws.Range("N4").Value = 10842468.31879392
$ws.Range("O4").Value = 10434123.96605083
$ws.Range("P4").Value = 0.01681289542120121
$ws.Range("Q4").Value = 0.01747087599732293

# Row 5 (B and C unchanged)
$ws.Range("D5").Value = 91
$ws.Range("E5").Value = 0.8666666666666667
$ws.Range("F5").Value = 0.8666666666666667
$ws.Range("G5").Value = 0.09511055620966978
$ws.Range("H5").Value = 0.08242914871504717
$ws.Range("I5").Value = 526523.2199554271
$ws.Range("J5").Value = 188388.8538193516
$ws.Range("L5").Value = 188388.8538193516
$ws.Range("M5").Value = 714912.0737747787
$ws.Range("N5").Value = 11280831.86545774
$ws.Range("O5").Value = 10870137.18213236
$ws.Range("P5").Value = 0.01669990795591983
$ws.Range("Q5").Value = 0.01733086258828576

# Row 6 (B and C unchanged)
$ws.Range("D6").Value = 92
$ws.Range("E6").Value = 0.8679245283018868
$ws.Range("F6").Value = 0.8679245283018868
$ws.Range("G6").Value = 0.09510156385804072
$ws.Range("H6").Value = 0.08254097995226177
$ws.Range("I6").Value = 551974.7101267558
$ws.Range("J6").Value = 198053.4162202651
$ws.Range("L6").Value = 198053.4162202651
$ws.Range("M6").Value = 750028.1263470209
$ws.Range("N6").Value = 11729815.62832147
$ws.Range("O6").Value = 11315350.10449633
$ws.Range("P6").Value = 0.01688461460059679
$ws.Range("Q6").Value = 0.01750307453072668
